# Updates KCOR/CI_Lower/CI_Upper values for specific (YearOfBirth) rows
# across all 7 sheets of KCOR_summary.xlsx, per commit 'this is awesome. We are done.'
# Values are written with a leading apostrophe so Excel keeps them as text
# (matching the original inlineStr/text storage of these cells).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2021_13")
$ws.Range("C13").Value = "'2.7164"
$ws.Range("D13").Value = "'1.823"
$ws.Range("E13").Value = "'4.048"

$ws.Range("C14").Value = "'17.0855"
$ws.Range("D14").Value = "'9.376"
$ws.Range("E14").Value = "'31.133"

$ws.Range("C27").Value = "'8.0417"
$ws.Range("D27").Value = "'2.583"
$ws.Range("E27").Value = "'25.039"

$ws.Range("C39").Value = "'1.0213"
$ws.Range("D39").Value = "'0.609"
$ws.Range("E39").Value = "'1.711"

$ws.Range("C40").Value = "'0.4707"
$ws.Range("D40").Value = "'0.131"
$ws.Range("E40").Value = "'1.686"

$ws = $wb.Worksheets.Item("2021_20")
$ws.Range("C13").Value = "'1.9442"
$ws.Range("D13").Value = "'1.418"
$ws.Range("E13").Value = "'2.666"

$ws.Range("C14").Value = "'10.4248"
$ws.Range("D14").Value = "'6.649"
$ws.Range("E14").Value = "'16.344"

$ws.Range("C27").Value = "'8.5882"
$ws.Range("D27").Value = "'5.036"
$ws.Range("E27").Value = "'14.645"

$ws.Range("C39").Value = "'1.4483"
$ws.Range("D39").Value = "'0.978"
$ws.Range("E39").Value = "'2.146"

$ws.Range("C40").Value = "'0.8238"
$ws.Range("D40").Value = "'0.416"
$ws.Range("E40").Value = "'1.630"

$ws = $wb.Worksheets.Item("2022_06")
$ws.Range("C13").Value = "'1.5388"
$ws.Range("D13").Value = "'0.880"
$ws.Range("E13").Value = "'2.690"

$ws.Range("C39").Value = "'0.2939"
$ws.Range("D39").Value = "'0.167"
$ws.Range("E39").Value = "'0.519"

$ws.Range("C65").Value = "'0.5775"
$ws.Range("D65").Value = "'0.321"
$ws.Range("E65").Value = "'1.038"

$ws = $wb.Worksheets.Item("2022_26")
$ws.Range("C14").Value = "'0.7206"
$ws.Range("D14").Value = "'0.369"
$ws.Range("E14").Value = "'1.405"

$ws.Range("C40").Value = "'0.8083"
$ws.Range("D40").Value = "'0.413"
$ws.Range("E40").Value = "'1.580"

$ws.Range("C66").Value = "'0.9365"
$ws.Range("D66").Value = "'0.469"
$ws.Range("E66").Value = "'1.869"

$ws = $wb.Worksheets.Item("2022_47")
$ws.Range("C13").Value = "'0.8948"
$ws.Range("D13").Value = "'0.475"
$ws.Range("E13").Value = "'1.685"

$ws.Range("C14").Value = "'1.4727"
$ws.Range("D14").Value = "'0.601"
$ws.Range("E14").Value = "'3.609"

$ws.Range("C39").Value = "'0.6503"
$ws.Range("D39").Value = "'0.342"
$ws.Range("E39").Value = "'1.236"

$ws.Range("C40").Value = "'0.6018"
$ws.Range("D40").Value = "'0.244"
$ws.Range("E40").Value = "'1.486"

$ws.Range("C65").Value = "'0.5430"
$ws.Range("D65").Value = "'0.284"
$ws.Range("E65").Value = "'1.040"

$ws.Range("C66").Value = "'0.6597"
$ws.Range("D66").Value = "'0.260"
$ws.Range("E66").Value = "'1.675"

$ws.Range("C92").Value = "'4.1487"
$ws.Range("D92").Value = "'1.530"
$ws.Range("E92").Value = "'11.251"

$ws.Range("C104").Value = "'1.7815"
$ws.Range("D104").Value = "'0.735"
$ws.Range("E104").Value = "'4.319"

$ws.Range("C118").Value = "'4.6807"
$ws.Range("D118").Value = "'1.715"
$ws.Range("E118").Value = "'12.776"

$ws.Range("C131").Value = "'4.2705"
$ws.Range("D131").Value = "'1.525"
$ws.Range("E131").Value = "'11.959"

